# Update the "想去人数" (want-to-go count) figures that changed between
# the previous gh-pages data snapshot and the newly generated one.
#
# Sheet "展览" (rId1 / sheet1.xml):
#   F3  2150 -> 2151   (苏州·燃梦)
#   F5  11317 -> 11318 (【会员购严选】苏州·端阳嘉年华动漫国潮文化节)
#   F9  11255 -> 11260 (苏州·萤火国潮文化节动漫品牌博览会)
#   F12 62 -> 63       (常熟·动漫游戏嘉年华)
#   F16 3465 -> 3466   (苏州·第十三届理想乡动漫展同人创作者大会)
#
# Sheet "全部类型" (rId4 / sheet4.xml) mirrors the same events at
# different row offsets:
#   F3  2150 -> 2151
#   F7  11317 -> 11318
#   F11 11255 -> 11260
#   F14 62 -> 63
#   F18 3465 -> 3466

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F3").Value  = 2151
$wsExhibition.Range("F5").Value  = 11318
$wsExhibition.Range("F9").Value  = 11260
$wsExhibition.Range("F12").Value = 63
$wsExhibition.Range("F16").Value = 3466

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value  = 2151
$wsAll.Range("F7").Value  = 11318
$wsAll.Range("F11").Value = 11260
$wsAll.Range("F14").Value = 63
$wsAll.Range("F18").Value = 3466
